$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these Price/Volume columns remain text so values like "1.00" or
# "0.999" are not coerced into numbers, matching the original inlineStr
# (text) cell contents.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "65.021.14"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "3.157.48"
$ws.Range("E3").Value = "  +3.10%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "574.17"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").Value = "150.06"
$ws.Range("E6").Value = "  +5.23%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.155.87"
$ws.Range("E8").Value = "  +3.12%  "
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("E10").Value = "  +4.20%  "
$ws.Range("D11").Value = "6.09"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "0.499"
$ws.Range("E12").Value = "  +3.66%  "
$ws.Range("D13").Value = "0.0000265"
$ws.Range("E13").Value = "  +14.28%  "
$ws.Range("D14").Value = "37.22"
$ws.Range("E14").Value = "  +5.16%  "
$ws.Range("D15").Value = "3.681.65"
$ws.Range("E15").Value = "  +3.29%  "
$ws.Range("D16").Value = "65.116.50"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "3.158.80"
$ws.Range("E17").Value = "  +3.46%  "
$ws.Range("D18").Value = "7.11"
$ws.Range("E18").Value = "  +4.36%  "
$ws.Range("E19").Value = "  +1.07%  "
$ws.Range("D20").Value = "507.26"
$ws.Range("E20").Value = "  +4.01%  "
$ws.Range("D21").Value = "14.92"
$ws.Range("E21").Value = "  +3.91%  "
$ws.Range("D22").Value = "0.719"
$ws.Range("E22").Value = "  +3.80%  "
$ws.Range("D23").Value = "15.33"
$ws.Range("E23").Value = "  +3.74%  "
$ws.Range("D24").Value = "7.73"
$ws.Range("E24").Value = "  +2.81%  "
$ws.Range("D25").Value = "84.40"
$ws.Range("E25").Value = "  +2.24%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "8.92"
$ws.Range("E27").Value = "  +8.64%  "
$ws.Range("E28").Value = "  +3.49%  "
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  +5.55%  "
$ws.Range("D30").Value = "2.80"
$ws.Range("E30").Value = "  +8.77%  "
$ws.Range("D31").Value = "27.65"
$ws.Range("E31").Value = "  +4.21%  "
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").Value = "  +3.21%  "
$ws.Range("D34").Value = "6.18"
$ws.Range("E34").Value = "  +7.08%  "
$ws.Range("D35").Value = "6.52"
$ws.Range("E35").Value = "  +3.81%  "
$ws.Range("D36").Value = "54.76"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "0.0902"
$ws.Range("E37").Value = "  +10.24%  "
$ws.Range("D38").Value = "466.17"
$ws.Range("E38").Value = "  +5.26%  "
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").Value = "3.00"
$ws.Range("E40").Value = "  +9.74%  "
$ws.Range("D41").Value = "8.68"
$ws.Range("E41").Value = "  +3.73%  "
$ws.Range("D42").Value = "3.059.64"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").Value = "2.44"
$ws.Range("E44").Value = "  +8.39%  "
$ws.Range("D45").Value = "0.283"
$ws.Range("E45").Value = "  +2.24%  "
$ws.Range("D46").Value = "28.75"
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("D47").Value = "{0}{1}{2}" -f "0.0", [char]0x2083, "0588"
$ws.Range("E47").Value = "  +13.19%  "
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("E50").Value = "  +5.15%  "
$ws.Range("D51").Value = "119.52"
$ws.Range("E51").Value = "  +1.77%  "
